$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style of an existing styled cell in column A down to the newly added rows (17-25)
$ws.Range("A16").Copy()
$ws.Range("A17:A25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Write the new/reordered data for rows 1-25 (col A = rank number, col B = UPN string)
$ws.Range("B1").Value = "upn"
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "R931258916027"
$ws.Range("A3").Value = 43
$ws.Range("B3").Value = "B931252109011"
$ws.Range("A4").Value = 54
$ws.Range("B4").Value = "Q931252909054"
$ws.Range("A5").Value = 47
$ws.Range("B5").Value = "Q931101109046"
$ws.Range("A6").Value = 41
$ws.Range("B6").Value = "L931412020028"
$ws.Range("A7").Value = 38
$ws.Range("B7").Value = "R931412017031"
$ws.Range("A8").Value = 37
$ws.Range("B8").Value = "V333218013124"
$ws.Range("A9").Value = 55
$ws.Range("B9").Value = "M931321009023"
$ws.Range("A10").Value = 4
$ws.Range("B10").Value = "D931252714054"
$ws.Range("A11").Value = 8
$ws.Range("B11").Value = "C931101008023"
$ws.Range("A12").Value = 14
$ws.Range("B12").Value = "E931252916073"
$ws.Range("A13").Value = 46
$ws.Range("B13").Value = "Y931325210014"
$ws.Range("A14").Value = 40
$ws.Range("B14").Value = "C373236909082"
$ws.Range("A15").Value = 17
$ws.Range("B15").Value = "V931414517045"
$ws.Range("A16").Value = 10
$ws.Range("B16").Value = "C931316110004"
$ws.Range("A17").Value = 11
$ws.Range("B17").Value = "U931258914007"
$ws.Range("A18").Value = 30
$ws.Range("B18").Value = "T887690719015"
$ws.Range("A19").Value = 39
$ws.Range("B19").Value = "T931252911047"
$ws.Range("A20").Value = 9
$ws.Range("B20").Value = "M931100509004"
$ws.Range("A21").Value = 42
$ws.Range("B21").Value = "W931321009038"
$ws.Range("A22").Value = 23
$ws.Range("B22").Value = "Q931252907052"
$ws.Range("A23").Value = 22
$ws.Range("B23").Value = "T931252108049"
$ws.Range("A24").Value = 52
$ws.Range("B24").Value = "K931383410019"
$ws.Range("A25").Value = 53
$ws.Range("B25").Value = "J931101109013"

$wb.Save()
